# Updated cryptos list on Fri May 26 15:23:13 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# on the active sheet. Price/volume cells hold plain text in the source
# data (e.g. "0.7460", "  +2.10%  "), so for any new Price value that
# looks like a bare number we first force the cell to Text format ("@")
# before assigning it - otherwise Excel would silently reinterpret it as
# a numeric value and drop meaningful trailing zeros (e.g. "0.7460" ->
# 0.746). Volume cells keep their padded "  +x.xx%  " shape, which Excel
# already treats as text, so no extra formatting is required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.985.56"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").Value = "1.847.17"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.36"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4673"
$ws.Range("E7").Value = "  +3.29%  "

$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07174"
$ws.Range("E9").Value = "  +1.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9305"
$ws.Range("E10").Value = "  +4.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.59"
$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07668"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").Value = "1.864.48"
$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.294"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.414"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.16"
$ws.Range("E16").Value = "  +2.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008588"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").Value = "27.009.93"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.39"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.036"
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +1.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.931"
$ws.Range("E24").Value = "  -1.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.66"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.04"
$ws.Range("E26").Value = "  +1.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.028"
$ws.Range("E27").Value = "  -1.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.06"
$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.931"
$ws.Range("E29").Value = "  +1.62%  "

$ws.Range("E30").Value = "  +1.88%  "

$ws.Range("E31").Value = "  +2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.853"
$ws.Range("E32").Value = "  +0.79%  "

$ws.Range("E33").Value = "  +6.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7460"
$ws.Range("E34").Value = "  +3.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.474"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.090"
$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.993"
$ws.Range("E37").Value = "  +3.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01939"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("E39").Value = "  +1.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5132"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.898"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1513"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.177"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.46"
$ws.Range("E44").Value = "  +4.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4709"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.12"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.607"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06048"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.33"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.10"
$ws.Range("E51").Value = "  +0.05%  "
